$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.991.00"
$ws.Range("E2").Value = "'  -0.91%  "
$ws.Range("D3").Value = "'1.574.17"
$ws.Range("E3").Value = "'  -2.06%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E5").Value = "'  +0.18%  "
$ws.Range("D6").Value = "'298.87"
$ws.Range("E6").Value = "'  -1.31%  "
$ws.Range("D7").Value = "'0.3735"
$ws.Range("E7").Value = "'  -0.85%  "
$ws.Range("D8").Value = "'0.3543"
$ws.Range("E8").Value = "'  -3.00%  "
$ws.Range("D9").Value = "'49.91"
$ws.Range("E9").Value = "'  +2.40%  "
$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = "'  +0.13%  "
$ws.Range("D11").Value = "'1.207"
$ws.Range("E11").Value = "'  -5.00%  "
$ws.Range("D12").Value = "'0.07926"
$ws.Range("E12").Value = "'  -1.94%  "
$ws.Range("D13").Value = "'21.72"
$ws.Range("E13").Value = "'  -5.81%  "
$ws.Range("D14").Value = "'6.387"
$ws.Range("E14").Value = "'  -2.85%  "
$ws.Range("D15").Value = "'7.235"
$ws.Range("E15").Value = "'  -4.44%  "
$ws.Range("D16").Value = "'0.00001215"
$ws.Range("E16").Value = "'  -4.17%  "
$ws.Range("D17").Value = "'1.578.86"
$ws.Range("E17").Value = "'  -1.93%  "
$ws.Range("D18").Value = "'91.47"
$ws.Range("E18").Value = "'  -0.18%  "
$ws.Range("D19").Value = "'0.06734"
$ws.Range("E19").Value = "'  -0.68%  "
$ws.Range("D20").Value = "'17.61"
$ws.Range("E20").Value = "'  -4.03%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "'  +0.12%  "
$ws.Range("D22").Value = "'6.334"
$ws.Range("E22").Value = "'  -3.79%  "
$ws.Range("D23").Value = "'22.992.17"
$ws.Range("E23").Value = "'  -0.95%  "
$ws.Range("D24").Value = "'12.56"
$ws.Range("E24").Value = "'  -4.10%  "
$ws.Range("D25").Value = "'2.358"
$ws.Range("E25").Value = "'  +0.19%  "
$ws.Range("D26").Value = "'2.808"
$ws.Range("E26").Value = "'  -3.75%  "
$ws.Range("D27").Value = "'20.49"
$ws.Range("E27").Value = "'  -2.94%  "
$ws.Range("D28").Value = "'146.85"
$ws.Range("E28").Value = "'  -2.44%  "
$ws.Range("D29").Value = "'5.155"
$ws.Range("E29").Value = "'  -1.88%  "
$ws.Range("D30").Value = "'131.06"
$ws.Range("E30").Value = "'  -1.11%  "
$ws.Range("E31").Value = "'  -3.39%  "
$ws.Range("D32").Value = "'6.494"
$ws.Range("E32").Value = "'  -6.80%  "
$ws.Range("D33").Value = "'1.752.04"
$ws.Range("E33").Value = "'  -2.03%  "
$ws.Range("D34").Value = "'0.9251"
$ws.Range("E34").Value = "'  -5.55%  "
$ws.Range("D35").Value = "'0.07287"
$ws.Range("E35").Value = "'  -5.84%  "
$ws.Range("D36").Value = "'0.02643"
$ws.Range("D37").Value = "'0.08714"
$ws.Range("E37").Value = "'  -1.72%  "
$ws.Range("D38").Value = "'9.863"
$ws.Range("E38").Value = "'  -2.35%  "
$ws.Range("D39").Value = "'0.2451"
$ws.Range("E39").Value = "'  -4.05%  "
$ws.Range("D40").Value = "'5.945"
$ws.Range("E40").Value = "'  -5.47%  "
$ws.Range("D41").Value = "'1.335"
$ws.Range("E41").Value = "'  -4.21%  "
$ws.Range("D42").Value = "'0.6822"
$ws.Range("E42").Value = "'  -4.72%  "
$ws.Range("D43").Value = "'11.74"
$ws.Range("E43").Value = "'  -8.32%  "
$ws.Range("D44").Value = "'14.51"
$ws.Range("E44").Value = "'  -8.14%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6287"
$ws.Range("E45").Value = "'  -4.87%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.957"
$ws.Range("E46").Value = "'  -0.68%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.231"
$ws.Range("E47").Value = "'  -3.12%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'129.76"
$ws.Range("E48").Value = "'  -1.08%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.07838"
$ws.Range("E49").Value = "'  -2.22%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "'1.183"
$ws.Range("E50").Value = "'  +0.90%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.158"
$ws.Range("E51").Value = "'  -1.57%  "
